{"js": "// Locate the \"LOM3205: Eletromagnetismo (Requisito)\" paragraph (last\n// requisite line). The three paragraphs that follow it -- a blank\n// paragraph, the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, and\n// the \"(c) 2020 ...\" footer line -- are removed, leaving only the\n// trailing blank paragraph that precedes the final page-break paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"LOM3205: Eletromagnetismo (Requisito)\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex >= 0 && targetIndex + 3 < paragraphs.items.length) {\n  const t1 = paragraphs.items[targetIndex + 1].text;\n  const t2 = paragraphs.items[targetIndex + 2].text;\n  const t3 = paragraphs.items[targetIndex + 3].text;\n\n  const matches =\n    t1 === \"\" &&\n    t2 === \"Ver no Jupiter Salvar em pdf Salvar em docx\" &&\n    t3.indexOf(\"Contact: luizeleno@usp.br\") !== -1;\n\n  if (matches) {\n    // Delete from the end backwards so earlier indices stay valid.\n    paragraphs.items[targetIndex + 3].delete();\n    paragraphs.items[targetIndex + 2].delete();\n    paragraphs.items[targetIndex + 1].delete();\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOM3205: Eletromagnetismo (Requisito)\" paragraph (last\n# requisite line). The three paragraphs that follow it -- a blank\n# paragraph, the \"Ver no Jupiter Salvar em pdf Salvar em docx\" line, and\n# the \"(c) 2020 ...\" footer line -- are removed, leaving only the\n# trailing blank paragraph that precedes the final page-break paragraph.\n$paras = $d.Paragraphs\n$targetIndex = -1\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"LOM3205: Eletromagnetismo (Requisito)\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -gt 0 -and ($targetIndex + 3) -le $paras.Count) {\n    $t1 = $paras.Item($targetIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n    $t2 = $paras.Item($targetIndex + 2).Range.Text.TrimEnd([char]13, [char]7)\n    $t3 = $paras.Item($targetIndex + 3).Range.Text.TrimEnd([char]13, [char]7)\n\n    $matches = ($t1 -eq \"\") -and\n               ($t2 -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") -and\n               ($t3 -like \"*Contact: luizeleno@usp.br*\")\n\n    if ($matches) {\n        $delStart = $paras.Item($targetIndex + 1).Range.Start\n        $delEnd = $paras.Item($targetIndex + 3).Range.End\n        $delRange = $d.Range($delStart, $delEnd)\n        $delRange.Delete()\n    }\n}\n"}
